# Auto-generated Excel COM-interop script to apply Yojimbo_Profits profit recalculation updates
# across multiple worksheets (ALC, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 1290.625
$ws.Range("I52").Value = 1625
$ws.Range("J52").Value = 1179.1666
$ws.Range("K52").Value = 4875
$ws.Range("L52").Value = 3537.4998
$ws.Range("M52").Value = -4715
$ws.Range("N52").Value = -3857.4998
$ws.Range("H69").Value = 3072.5
$ws.Range("J69").Value = 2800
$ws.Range("L69").Value = 8400
$ws.Range("N69").Value = -10148
$ws.Range("H72").Value = 3072.5
$ws.Range("J72").Value = 2800
$ws.Range("L72").Value = 25200
$ws.Range("N72").Value = -33936
$ws.Range("H103").Value = 1376.125
$ws.Range("I103").Value = 1368
$ws.Range("J103").Value = 1381
$ws.Range("K103").Value = 4104
$ws.Range("L103").Value = 4143
$ws.Range("M103").Value = -3518
$ws.Range("N103").Value = -5315
$ws.Range("H125").Value = 1708.6364
$ws.Range("I125").Value = 1844.6
$ws.Range("J125").Value = 1595.3334
$ws.Range("K125").Value = 16601.4
$ws.Range("L125").Value = 14358.0006
$ws.Range("M125").Value = -14141.4
$ws.Range("N125").Value = -19278.0006
$ws.Range("H137").Value = 3020.238
$ws.Range("I137").Value = 2816.9697
$ws.Range("J137").Value = 3765.5557
$ws.Range("K137").Value = 8450.9091
$ws.Range("L137").Value = 11296.6671
$ws.Range("M137").Value = -5900.909100000001
$ws.Range("N137").Value = -16396.6671

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 992.2308
$ws.Range("I94").Value = 734.32355
$ws.Range("K94").Value = 734.32355
$ws.Range("M94").Value = -283.32355

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23320.451
$ws.Range("I31").Value = 47030.457
$ws.Range("J31").Value = 2244.889
$ws.Range("K31").Value = 47030.457
$ws.Range("L31").Value = 2244.889
$ws.Range("M31").Value = -46735.457
$ws.Range("N31").Value = -2834.889
$ws.Range("H34").Value = 23320.451
$ws.Range("I34").Value = 47030.457
$ws.Range("J34").Value = 2244.889
$ws.Range("K34").Value = 47030.457
$ws.Range("L34").Value = 2244.889
$ws.Range("M34").Value = -46828.457
$ws.Range("N34").Value = -2648.889
$ws.Range("H99").Value = 1831.091
$ws.Range("I99").Value = 1733.3334
$ws.Range("J99").Value = 1948.4
$ws.Range("K99").Value = 1733.3334
$ws.Range("L99").Value = 1948.4
$ws.Range("M99").Value = -235.3334
$ws.Range("N99").Value = -4944.4
$ws.Range("H126").Value = 1831.091
$ws.Range("I126").Value = 1733.3334
$ws.Range("J126").Value = 1948.4
$ws.Range("K126").Value = 5200.0002
$ws.Range("L126").Value = 5845.200000000001
$ws.Range("M126").Value = -2730.0002
$ws.Range("N126").Value = -10785.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 12059
$ws.Range("J62").Value = 12059
$ws.Range("L62").Value = 36177
$ws.Range("N62").Value = -37549
$ws.Range("H63").Value = 4000
$ws.Range("J63").Value = 5000
$ws.Range("L63").Value = 15000
$ws.Range("N63").Value = -16498
$ws.Range("H64").Value = 11632.167
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 11632.167
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 34896.501
$ws.Range("M64").ClearContents() | Out-Null
$ws.Range("N64").Value = -35436.501
$ws.Range("H65").Value = 12059
$ws.Range("J65").Value = 12059
$ws.Range("L65").Value = 108531
$ws.Range("N65").Value = -115395
$ws.Range("H66").Value = 4000
$ws.Range("J66").Value = 5000
$ws.Range("L66").Value = 45000
$ws.Range("N66").Value = -52488
$ws.Range("H67").Value = 11632.167
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 11632.167
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 34896.501
$ws.Range("M67").ClearContents() | Out-Null
$ws.Range("N67").Value = -36768.501
$ws.Range("H92").Value = 860.8
$ws.Range("I92").Value = 502
$ws.Range("J92").Value = 1100
$ws.Range("K92").Value = 1506
$ws.Range("L92").Value = 3300
$ws.Range("M92").Value = -258
$ws.Range("N92").Value = -5796
$ws.Range("H93").Value = 6500
$ws.Range("I93").Value = 3000
$ws.Range("J93").Value = 10000
$ws.Range("K93").Value = 9000
$ws.Range("L93").Value = 30000
$ws.Range("M93").Value = -7128
$ws.Range("N93").Value = -33744
$ws.Range("H94").Value = 2924
$ws.Range("I94").Value = 2888
$ws.Range("K94").Value = 8664
$ws.Range("M94").Value = -7988

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1169
$ws.Range("I126").Value = 1102.6666
$ws.Range("J126").Value = 1268.5
$ws.Range("K126").Value = 3307.9998
$ws.Range("L126").Value = 3805.5
$ws.Range("M126").Value = -837.9998
$ws.Range("N126").Value = -8745.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2652.3333
$ws.Range("I7").Value = 1352
$ws.Range("J7").Value = 3302.5
$ws.Range("K7").Value = 1352
$ws.Range("L7").Value = 3302.5
$ws.Range("M7").Value = -1240
$ws.Range("N7").Value = -3526.5
$ws.Range("H16").Value = 883.3333
$ws.Range("I16").Value = 826.4706
$ws.Range("J16").Value = 1850
$ws.Range("K16").Value = 826.4706
$ws.Range("L16").Value = 1850
$ws.Range("M16").Value = -656.4706
$ws.Range("N16").Value = -2190
$ws.Range("H40").Value = 5523.2354
$ws.Range("I40").Value = 5010.5557
$ws.Range("J40").Value = 6100
$ws.Range("K40").Value = 5010.5557
$ws.Range("L40").Value = 6100
$ws.Range("M40").Value = -4874.5557
$ws.Range("N40").Value = -6372
$ws.Range("H100").Value = 2767.3242
$ws.Range("I100").Value = 2239.1
$ws.Range("J100").Value = 2962.963
$ws.Range("K100").Value = 2239.1
$ws.Range("L100").Value = 2962.963
$ws.Range("M100").Value = -1698.1
$ws.Range("N100").Value = -4044.963
$ws.Range("H122").Value = 3633.8572
$ws.Range("I122").Value = 4188
$ws.Range("J122").Value = 2895
$ws.Range("K122").Value = 12564
$ws.Range("L122").Value = 8685
$ws.Range("M122").Value = -10114
$ws.Range("N122").Value = -13585
$ws.Range("H126").Value = 2652.3333
$ws.Range("I126").Value = 1352
$ws.Range("J126").Value = 3302.5
$ws.Range("K126").Value = 4056
$ws.Range("L126").Value = 9907.5
$ws.Range("M126").Value = -1586
$ws.Range("N126").Value = -14847.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5159975.5
$ws.Range("I122").Value = 6347478.5
$ws.Range("J122").Value = 4632196.5
$ws.Range("K122").Value = 19042435.5
$ws.Range("L122").Value = 13896589.5
$ws.Range("M122").Value = -19039985.5
$ws.Range("N122").Value = -13901489.5
